# Applies the cryptos.xlsx price/volume/date refresh described in the commit
# "Updated symbol list on Thu Dec 29 23:14:28 UTC 2022 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.13"
$ws.Range("G2").Value = "'23"
$ws.Range("D3").Value = "'24.20"
$ws.Range("G3").Value = "'23"
$ws.Range("D4").Value = "'5.277"
$ws.Range("G4").Value = "'23"
$ws.Range("D5").Value = "'0.05799"
$ws.Range("G5").Value = "'23"
$ws.Range("D6").Value = "'6.500"
$ws.Range("G6").Value = "'23"
$ws.Range("D7").Value = "'3.130"
$ws.Range("G7").Value = "'23"
$ws.Range("D8").Value = "'0.8168"
$ws.Range("G8").Value = "'23"
$ws.Range("D9").Value = "'0.8533"
$ws.Range("G9").Value = "'23"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.009804"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("G10").Value = "'23"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1359"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").Value = "'23"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06942"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'23"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03140"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'23"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02881"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "'23"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09395"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").Value = "'23"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.738"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").Value = "'23"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001517"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").Value = "'23"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04672"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "'23"
$ws.Range("D19").Value = "'0.006280"
$ws.Range("G19").Value = "'23"
$ws.Range("D20").Value = "'0.001235"
$ws.Range("G20").Value = "'23"
$ws.Range("D21").Value = "'0.004628"
$ws.Range("G21").Value = "'23"
$ws.Range("D22").Value = "'0.00006895"
$ws.Range("G22").Value = "'23"
$ws.Range("D23").Value = "'3.499"
$ws.Range("G23").Value = "'23"
$ws.Range("D24").Value = "'2.148"
$ws.Range("G24").Value = "'23"
$ws.Range("G25").Value = "'23"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("G26").Value = "'23"
$ws.Range("G27").Value = "'23"
$ws.Range("D28").Value = "'0.0002330"
$ws.Range("G28").Value = "'23"
$ws.Range("G29").Value = "'23"
$ws.Range("G30").Value = "'23"
$ws.Range("G31").Value = "'23"
$ws.Range("G32").Value = "'23"
$ws.Range("G33").Value = "'23"
$ws.Range("G34").Value = "'23"
$ws.Range("G35").Value = "'23"
$ws.Range("G36").Value = "'23"
$ws.Range("G37").Value = "'23"
$ws.Range("G38").Value = "'23"
$ws.Range("G39").Value = "'23"
$ws.Range("D40").Value = "'0.03659"
$ws.Range("G40").Value = "'23"
$ws.Range("D41").Value = "'0.006255"
$ws.Range("G41").Value = "'23"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("G42").Value = "'23"
$ws.Range("D43").Value = "'0.003397"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("G43").Value = "'23"
$ws.Range("D44").Value = "'0.007475"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("G44").Value = "'23"
$ws.Range("D45").Value = "'0.00005256"
$ws.Range("G45").Value = "'23"
$ws.Range("G46").Value = "'23"
$ws.Range("D47").Value = "'0.3697"
$ws.Range("G47").Value = "'23"
$ws.Range("G48").Value = "'23"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("G49").Value = "'23"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("G50").Value = "'23"
$ws.Range("G51").Value = "'23"
